$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.573.74"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "2.615.37"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Formula = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Formula = "'534.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Formula = "'142.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Formula = "'0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Formula = "'6.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("D10").Formula = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "3.078.07"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "58.514.96"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Formula = "'20.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "2.611.34"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Formula = "'0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Formula = "'4.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Formula = "'334.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Formula = "'10.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Formula = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Formula = "'66.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").Formula = "'0.418"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").Formula = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").Formula = "'7.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "0.0₃0733"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Formula = "'5.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Formula = "'153.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").Formula = "'18.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Formula = "'3.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Formula = "'0.843"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").Formula = "'1.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").Formula = "'0.814"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Formula = "'3.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Formula = "'282.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").Formula = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Formula = "'0.594"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Formula = "'0.0945"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Formula = "'18.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").Value = "1.942.42"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Formula = "'4.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Formula = "'17.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("D51").Formula = "'113.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
